# Insert 4 new data rows at the top of the data block (rows 63-66),
# pushing all subsequent rows down by 4 (dimension grows from T137 to T141).
# Then populate the 4 new rows with the new day's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63..137 down to 67..141 by inserting 4 blank rows at 63.
$ws.Range("A63:T66").Insert()

# Row 63: Especial
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 45219
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100112025
$ws.Range("J63").Value = "Frutilla"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Especial"
$ws.Range("M63").Value = 280
$ws.Range("N63").Value = 8000
$ws.Range("O63").Value = 9000
$ws.Range("P63").Value = 8643
$ws.Range("Q63").Value = "$/bandeja 3 kilos"
$ws.Range("R63").Value = "Región de Arica y Parinacota"
$ws.Range("S63").Value = 2881
$ws.Range("T63").Value = 3

# Row 64: Primera
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 45219
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100101
$ws.Range("H64").Value = "Berries"
$ws.Range("I64").Value = 100112025
$ws.Range("J64").Value = "Frutilla"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 450
$ws.Range("N64").Value = 6000
$ws.Range("O64").Value = 7000
$ws.Range("P64").Value = 6556
$ws.Range("Q64").Value = "$/bandeja 3 kilos"
$ws.Range("R64").Value = "Región de Arica y Parinacota"
$ws.Range("S64").Value = 2185
$ws.Range("T64").Value = 3

# Row 65: Segunda
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 45219
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100112025
$ws.Range("J65").Value = "Frutilla"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 550
$ws.Range("N65").Value = 4000
$ws.Range("O65").Value = 5000
$ws.Range("P65").Value = 4545
$ws.Range("Q65").Value = "$/bandeja 3 kilos"
$ws.Range("R65").Value = "Región de Arica y Parinacota"
$ws.Range("S65").Value = 1515
$ws.Range("T65").Value = 3

# Row 66: Tercera
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 45219
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100101
$ws.Range("H66").Value = "Berries"
$ws.Range("I66").Value = 100112025
$ws.Range("J66").Value = "Frutilla"
$ws.Range("K66").Value = "Sin especificar"
$ws.Range("L66").Value = "Tercera"
$ws.Range("M66").Value = 200
$ws.Range("N66").Value = 2000
$ws.Range("O66").Value = 3000
$ws.Range("P66").Value = 2400
$ws.Range("Q66").Value = "$/bandeja 3 kilos"
$ws.Range("R66").Value = "Región de Arica y Parinacota"
$ws.Range("S66").Value = 800
$ws.Range("T66").Value = 3
